$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "time_taken", styled like the other header cells
# (bold font, thin border all around, centered horizontally, top vertically)
$ws.Range("F1").Value = "time_taken"
$f1 = $ws.Range("F1")
$f1.Font.Bold = $true
$f1.HorizontalAlignment = -4108
$f1.VerticalAlignment = -4160
$f1.Borders.LineStyle = 1
$f1.Borders.Weight = 2

# Timestamp values for data rows 2-34 (plain, unstyled cells), matching the source diff
$timeTaken = @(
    "2021-10-05 13:42:05.800437",
    "2021-10-05 13:42:05.800450",
    "2021-10-05 13:42:05.800454",
    "2021-10-05 13:42:05.800457",
    "2021-10-05 13:42:05.800461",
    "2021-10-05 13:42:05.800464",
    "2021-10-05 13:42:05.800467",
    "2021-10-05 13:42:05.800470",
    "2021-10-05 13:42:05.800473",
    "2021-10-05 13:42:05.800476",
    "2021-10-05 13:42:05.800479",
    "2021-10-05 13:42:05.800482",
    "2021-10-05 13:42:05.800485",
    "2021-10-05 13:42:05.800488",
    "2021-10-05 13:42:05.800491",
    "2021-10-05 13:42:05.800494",
    "2021-10-05 13:42:05.800498",
    "2021-10-05 13:42:05.800501",
    "2021-10-05 13:42:05.800504",
    "2021-10-05 13:42:05.800507",
    "2021-10-05 13:42:05.800510",
    "2021-10-05 13:42:05.800513",
    "2021-10-05 13:42:05.800515",
    "2021-10-05 13:42:05.800518",
    "2021-10-05 13:42:05.800522",
    "2021-10-05 13:42:05.800525",
    "2021-10-05 13:42:05.800528",
    "2021-10-05 13:42:05.800531",
    "2021-10-05 13:42:05.800534",
    "2021-10-05 13:42:05.800537",
    "2021-10-05 13:42:05.800540",
    "2021-10-05 13:42:05.800543",
    "2021-10-05 13:42:05.800546"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
